$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as plain text in the source data
# (e.g. thousands-dot formatted numbers like "47.195.99"). For any new
# price that Excel would otherwise auto-parse as a number, force the cell
# to Text format first so the literal string is preserved.
$ws.Range("D2").Value = "47.195.99"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "2.483.00"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.14"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.54"
$ws.Range("E6").Value = "  +1.87%  "
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.40"
$ws.Range("E10").Value = "  +3.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0806"
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.22"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("D15").Value = "2.873.95"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").Value = "2.479.64"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").Value = "47.112.29"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.72"
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.60"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").Value = "0.0₃0928"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("E22").Value = "  +12.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.20"
$ws.Range("E23").Value = "  -1.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "244.67"
$ws.Range("E24").Value = "  -3.09%  "
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  -3.57%  "
$ws.Range("E28").Value = "  +3.33%  "
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.136"
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.24"
$ws.Range("E31").Value = "  -2.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.38"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.15"
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0775"
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.62"
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("E39").Value = "  -1.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.51"
$ws.Range("E40").Value = "  +5.43%  "
$ws.Range("E41").Value = "  -1.15%  "
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "118.25"
$ws.Range("E43").Value = "  -4.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0294"
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("D45").Value = "1.981.41"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("E47").Value = "  -6.85%  "
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.10"
$ws.Range("E50").Value = "  -5.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.40"
$ws.Range("E51").Value = "  +2.72%  "
